# Generate Report for Handoff
#
# Localization status moved from "In Translation" to "Ready for handoff":
#   - Update the status cells on the Overview sheet (per-language status)
#     and on each language sheet's Status column.
#   - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps to the new handoff generation time.
#   - Widen the status columns a bit so the longer "Ready for handoff"
#     label fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet: E2 (zh-cn status), F2 (de-de status), G2 (generate date) ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-30 00:40:52"

# --- zh-cn sheet: C2 (status), H2 (latest handoff datetime) ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-30 00:40:48"

# --- de-de sheet: C2 (status), H2 (latest handoff datetime) ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-30 00:40:52"

# --- Widen the status columns to fit the longer "Ready for handoff" text ---
# (Excel's ColumnWidth setter snaps to whole-pixel character widths, so the
# input width is chosen to land on the closest achievable value.)
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336
$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
